$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$title = $s.Shapes.Item(1)
$tr = $title.TextFrame.TextRange

$tr.Characters(1, 4).Font.Italic = $true
